# Add a new paragraph "Mani" right after the "Karthick" paragraph
# (and before the "Iam karthick" paragraph).
$d = $word.ActiveDocument

# The "Karthick" paragraph is the 2nd paragraph in the document
# ("Vikram", "Karthick", "Iam karthick").
$karthickPara = $d.Paragraphs.Item(2)

# Insert a brand-new empty paragraph right after it, splitting off a
# fresh paragraph mark that inherits the same paragraph/run formatting
# (en-US language), matching the source document's pattern.
$karthickPara.Range.InsertParagraphAfter()

# Re-fetch the document/paragraph collection since the collection was
# mutated, then fill in the text of the newly created paragraph.
$newPara = $word.ActiveDocument.Paragraphs.Item(3)
$newPara.Range.Text = "Mani"

Write-Output "Inserted 'Mani' paragraph after 'Karthick'"
